# Update the pre-approved automotive credit template:
# wrap the comma-separated list of commercial-house RUCs in square
# brackets so it is stored as a single array-like value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

$newValue = "[1002003004001,1002003004002]"

$ws.Range("V2").Value = $newValue
$ws.Range("V3").Value = $newValue

# Move/leave the active selection where the author left it while editing.
$ws.Range("H17").Select() | Out-Null
